$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I11").Value = 'sd'
$ws.Range("J11").Value = 'Statement-non-opinion'
$ws.Range("I18").Value = 'ba'
$ws.Range("J18").Value = 'Appreciation'
$ws.Range("I21").Value = 'ba'
$ws.Range("J21").Value = 'Appreciation'
$ws.Range("I55").Value = 'aa'
$ws.Range("J55").Value = 'Agree/Accept'
$ws.Range("I74").Value = 'aa'
$ws.Range("J74").Value = 'Agree/Accept'
$ws.Range("I75").Value = 'sd'
$ws.Range("J75").Value = 'Statement-non-opinion'
$ws.Range("I80").Value = 'sv'
$ws.Range("J80").Value = 'Statement-opinion'
$ws.Range("I81").Value = 'sv'
$ws.Range("J81").Value = 'Statement-opinion'
$ws.Range("I93").Value = 'aa'
$ws.Range("J93").Value = 'Agree/Accept'
$ws.Range("I97").Value = 'sd'
$ws.Range("J97").Value = 'Statement-non-opinion'
$ws.Range("I101").Value = 'sv'
$ws.Range("J101").Value = 'Statement-opinion'
$ws.Range("I119").Value = 'ba'
$ws.Range("J119").Value = 'Appreciation'
$ws.Range("I126").Value = 'sd'
$ws.Range("J126").Value = 'Statement-non-opinion'
$ws.Range("I133").Value = 'ba'
$ws.Range("J133").Value = 'Appreciation'
$ws.Range("I135").Value = 'b'
$ws.Range("J135").Value = 'Acknowledge (Backchannel)'
$ws.Range("I136").Value = 'sd'
$ws.Range("J136").Value = 'Statement-non-opinion'
$ws.Range("I140").Value = 'sd'
$ws.Range("J140").Value = 'Statement-non-opinion'
$ws.Range("I150").Value = 'aa'
$ws.Range("J150").Value = 'Agree/Accept'
$ws.Range("I152").Value = 'sd'
$ws.Range("J152").Value = 'Statement-non-opinion'
$ws.Range("I164").Value = 'aa'
$ws.Range("J164").Value = 'Agree/Accept'
$ws.Range("I171").Value = 'aa'
$ws.Range("J171").Value = 'Agree/Accept'
$ws.Range("I175").Value = 'aa'
$ws.Range("J175").Value = 'Agree/Accept'
$ws.Range("I198").Value = 'sv'
$ws.Range("J198").Value = 'Statement-opinion'
$ws.Range("I207").Value = 'b'
$ws.Range("J207").Value = 'Acknowledge (Backchannel)'
$ws.Range("I209").Value = 'sv'
$ws.Range("J209").Value = 'Statement-opinion'
$ws.Range("I212").Value = '%'
$ws.Range("J212").Value = 'Uninterpretable'
$ws.Range("I226").Value = '%'
$ws.Range("J226").Value = 'Uninterpretable'
$ws.Range("I228").Value = 'sd'
$ws.Range("J228").Value = 'Statement-non-opinion'
$ws.Range("I236").Value = 'aa'
$ws.Range("J236").Value = 'Agree/Accept'
$ws.Range("I237").Value = 'sv'
$ws.Range("J237").Value = 'Statement-opinion'
$ws.Range("I239").Value = 'sv'
$ws.Range("J239").Value = 'Statement-opinion'
$ws.Range("I250").Value = 'aa'
$ws.Range("J250").Value = 'Agree/Accept'
$ws.Range("I256").Value = 'b'
$ws.Range("J256").Value = 'Acknowledge (Backchannel)'
$ws.Range("I272").Value = 'aa'
$ws.Range("J272").Value = 'Agree/Accept'
$ws.Range("I273").Value = 'sd'
$ws.Range("J273").Value = 'Statement-non-opinion'
$ws.Range("I274").Value = 'sd'
$ws.Range("J274").Value = 'Statement-non-opinion'
$ws.Range("I285").Value = 'sd'
$ws.Range("J285").Value = 'Statement-non-opinion'
